# Applies the attendance_reports sync update described in the commit:
# "Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-25 09:13:20"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: assign literal text that LOOKS like a number/percentage/fraction
# without letting Excel's input parser convert it into a real number (which
# would also silently swap the cell's style for a numeric-formatted one).
# We force-format a scratch cell as Text once, drop the literal string in
# it, then copy only the VALUE (not the format) onto the destination cell -
# so the destination keeps its original style untouched.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1000")
$scratch.NumberFormat = "@"

function Set-LiteralText {
    param($addr, $text)
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# Column I got narrower (14 -> 10 chars)
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 9.17

# ---------------------------------------------------------------------------
# "Recorded By" email-list reorders (same final set of recipients, new order)
# ---------------------------------------------------------------------------
$emailsA1 = "nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
foreach ($addr in @("G2", "G17", "G92", "G107")) {
    $ws.Range($addr).Value = $emailsA1
}

$emailsHisto = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
foreach ($addr in @("G7", "G22", "G112")) {
    $ws.Range($addr).Value = $emailsHisto
}

$emailsA3A4 = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
foreach ($addr in @("G32", "G47")) {
    $ws.Range($addr).Value = $emailsA3A4
}

$emailsB1B2 = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
foreach ($addr in @("G62", "G77")) {
    $ws.Range($addr).Value = $emailsB1B2
}

# ---------------------------------------------------------------------------
# Class Statistics block (K/L columns) numeric updates
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 14      # Recorded Sessions
$ws.Range("L7").Value = 0       # Missing Sessions
Set-LiteralText "L9"  "11.7%"   # Coverage %
Set-LiteralText "L10" "54.6%"   # Average Attendance %

# ---------------------------------------------------------------------------
# Per-group table row 16 (A2 / week 2) - S16 (Avg Attendance %) refreshed
# ---------------------------------------------------------------------------
Set-LiteralText "S16" "56.7%"

# ---------------------------------------------------------------------------
# Per-group table row 21 (B3) - Recorded/Missing counts + derived percentages
# ---------------------------------------------------------------------------
$ws.Range("O21").Value = 2
$ws.Range("P21").Value = 0
Set-LiteralText "R21" "13.3%"
Set-LiteralText "S21" "36.6%"

# ---------------------------------------------------------------------------
# Row 22 (A2 / B2 HISTOLOGY) - attendance count + derived percentage
# ---------------------------------------------------------------------------
$ws.Range("H22").Value = "117/217"
Set-LiteralText "S22" "43.8%"

# ---------------------------------------------------------------------------
# Row 97 (B3 HISTOLOGY) switched from "Not Recorded" (pink style) to
# "Recorded" (green style, same as the other already-recorded rows), and
# picked up recorders + an attendance count.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A97:I97").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G97").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("H97").Value = "1/224"
$ws.Range("I97").Value = "Recorded"

# ---------------------------------------------------------------------------
# Row 112 (B4 HISTOLOGY) - attendance count refreshed
# ---------------------------------------------------------------------------
$ws.Range("H112").Value = "97/226"

# Clean up the scratch cell used for forcing literal text.
$scratch.Clear()

$wb.Save()
